$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Seed the new event row's placeholder cells by copying from row 18 BEFORE
# row 18 itself is edited below:
#  - A18 already holds the card number "18" as text; reuse it verbatim.
#  - B18:K18 and N18 are still blank placeholders at this point, matching
#    what the freshly appended event row should start with.
$ws.Range("A18").Copy($ws.Range("A19"))
$ws.Range("B18:K18").Copy($ws.Range("B19:K19"))
$ws.Range("N18").Copy($ws.Range("N19"))

# Row 18 previously had blank placeholders in B..K and N; they now carry the
# literal "nan" marker used throughout the rest of the sheet.
$ws.Range("B18").Value = "nan"
$ws.Range("C18").Value = "nan"
$ws.Range("D18").Value = "nan"
$ws.Range("E18").Value = "nan"
$ws.Range("F18").Value = "nan"
$ws.Range("G18").Value = "nan"
$ws.Range("H18").Value = "nan"
$ws.Range("I18").Value = "nan"
$ws.Range("J18").Value = "nan"
$ws.Range("K18").Value = "nan"
$ws.Range("N18").Value = "nan"

# Append the new service event as row 19.
$ws.Range("L19").Value = "20\8\2025"
$ws.Range("M19").Value = "785 t"
$ws.Range("O19").Value = "تم تغيير جريده1  وجريده اليكران(90)"
$ws.Range("P19").Value = "الخبير"
